$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full branch/part/chapter lookup table (header + 47 data rows).
$data = @(
  @('branch', 'part', 'chapter'),
  @('01_whole-app-game', 'Intro', 'Whole app game'),
  @('02.1_shiny-app', 'Intro', 'Shiny'),
  @('02.2_movies-app', 'Intro', 'Shiny'),
  @('02.3_proj-app', 'Intro', 'Shiny'),
  @('03.1_description', 'Intro', 'Packages'),
  @('03.2_rproj', 'Intro', 'Packages'),
  @('03.3_create-package', 'Intro', 'Packages'),
  @('04_devtools', 'Intro', 'Development'),
  @('05_roxygen2', 'App-packages', 'Documentation'),
  @('06.1_pkg-exports', 'App-packages', 'Dependencies'),
  @('06.2_pkg-imports', 'App-packages', 'Dependencies'),
  @('07_data', 'App-packages', 'Data'),
  @('08_launch-app', 'App-packages', 'Launch'),
  @('09.1_inst-www', 'App-packages', 'External files'),
  @('09.2_inst-bslib', 'App-packages', 'External files'),
  @('09.3_inst-dev', 'App-packages', 'External files'),
  @('09.4_inst-prod', 'App-packages', 'External files'),
  @('11_tests-specs', 'Tests', 'Specifications'),
  @('12.1_tests-fixtures', 'Tests', 'Test tools'),
  @('12.2_tests-helpers', 'Tests', 'Test tools'),
  @('13_tests-modules', 'Tests', 'Testing modules'),
  @('14_tests-system', 'Tests', 'System tests'),
  @('15_docker', 'Deploy', 'Docker'),
  @('16.1_cicd-style', 'Deploy', 'CI/CD'),
  @('16.2_cicd-shiny', 'Deploy', 'CI/CD'),
  @('16.3_cicd-docker', 'Deploy', 'CI/CD'),
  @('17_golem', 'Frameworks', 'golem'),
  @('18_leprechaun', 'Frameworks', 'leprechaun'),
  @('19_rhino', 'Frameworks', 'rhino'),
  @('20_css', 'Non-R Code', 'CSS'),
  @('21_js', 'Non-R Code', 'JavaScript'),
  @('22_python', 'Non-R Code', 'Python'),
  @('23.1_debug-error', 'Special Topics', 'Debugging'),
  @('23.2_debug-selected_vars', 'Special Topics', 'Debugging'),
  @('23.3_debug-var_inputs', 'Special Topics', 'Debugging'),
  @('23.4_debug-scatter_plot', 'Special Topics', 'Debugging'),
  @('23.5_debug-print', 'Special Topics', 'Debugging'),
  @('24.1.0_reactive-values', 'Special Topics', 'Values vs. data'),
  @('24.1.1_step_01', 'Special Topics', 'Values vs. data'),
  @('24.1.2_step_02', 'Special Topics', 'Values vs. data'),
  @('24.1.3_step_03', 'Special Topics', 'Values vs. data'),
  @('24.1.4_step_04', 'Special Topics', 'Values vs. data'),
  @('24.2_user-data', 'Special Topics', 'Values vs. data'),
  @('24.2.0_user-data', 'Special Topics', 'Values vs. data'),
  @('24.2.1_step_01', 'Special Topics', 'Values vs. data'),
  @('24.2.2_step_02', 'Special Topics', 'Values vs. data'),
  @('A.E_mocks-snapshots', 'Appendix', 'Test mocks and snapshots')
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
  }
}

# Re-select the used range (now A1:C48) and widen the text columns so they
# fit the longer branch / chapter labels (matches the ~22.83-char bestFit
# width Excel computes for the new longest entries in columns A and C).
$ws.Range("A1:C48").Select()
$ws.Columns.Item(1).ColumnWidth = 22
$ws.Columns.Item(3).ColumnWidth = 22

$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
